$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing existing rows 5..84 down to 6..85
$ws.Range("A5:T5").Insert()

# Populate the newly inserted row 5 with the new record
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 44756
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100108
$ws.Cells.Item(5, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(5, 9).Value = 100108007
$ws.Cells.Item(5, 10).Value = "Coco"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 25
$ws.Cells.Item(5, 14).Value = 28000
$ws.Cells.Item(5, 15).Value = 28000
$ws.Cells.Item(5, 16).Value = 28000
$ws.Cells.Item(5, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(5, 18).Value = "Perú"
$ws.Cells.Item(5, 19).Value = 1400
$ws.Cells.Item(5, 20).Value = 20
